$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Working Papers")

$newTitle = "Risk Aversion and Barriers to Firm Growth: Experimental Evidence from Small Retailers"
$newAbstract = "Many enterprises in low and middle-income countries (LMICs) are owner-operated, so uncertain investments directly impact household consumption. This paper asks whether this link results in risk aversion that prevents firms from making uncertain investments needed to grow. I develop a model of small firm learning which shows that risk aversion can impede product adoption by preventing experimentation with new goods. I test the model using two field experiments in Kenya. The first directly tests whether risk aversion affects product stocking. Firms were offered new product inventory with an insurance contract that would be strictly dominated under risk neutrality. This increased adoption by 50\%, rejecting firm risk neutrality. The second experiment examines longer-run effects. Temporarily inducing firms to try selling a new product with a supplier returns policy led to a 70\% increase in stocking \textit{after} the intervention ended. The policy encouraged uncertain but optimistic firms to experiment, and they resolved uncertainty through learning. I validate the importance of learning by showing that firms were more willing to stock when the continuation value of learning increased, and when firms received information about demand. These results show that risk aversion can prevent small firms from making uncertain but high expected return investments, creating a barrier to growth that is not captured by standard models that assume risk neutral firms."

$ws.Range("C2").Value = $newAbstract
$ws.Range("A2").Value = $newTitle

$ws.Range("B2").Select() | Out-Null
